$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C values from 0.05 to 0.2 for rows 1-9
$ws.Range("C1:C9").Value = 0.2

# Update the selected cell/active cell to F7
$ws.Range("F7").Select()
